$d = $word.ActiveDocument

$replacements = @(
    @("N = 90,258", "N = 90,237"),
    @("2,017 (2.2)", "2,016 (2.2)"),
    @("1,284 (1.4)", "1,283 (1.4)"),
    @("2,680 (3.0)", "2,679 (3.0)"),
    @("87,578 (97)", "87,558 (97)"),
    @("51,733 (57)", "51,723 (57)"),
    @("38,525 (43)", "38,514 (43)"),
    @("7,120 (7.9)", "7,118 (7.9)"),
    @("22,117 (25)", "22,114 (25)"),
    @("21,321 (24)", "21,318 (24)"),
    @("39,700 (44)", "39,687 (44)"),
    @("11,586 (13)", "11,584 (13)"),
    @("19,663 (22)", "19,661 (22)"),
    @("23,668 (26)", "23,663 (26)"),
    @("20,901 (23)", "20,894 (23)"),
    @("6,074 (6.7)", "6,072 (6.7)"),
    @("8,366 (9.3)", "8,363 (9.3)"),
    @("52,105 (58)", "52,094 (58)"),
    @("32,075 (36)", "32,065 (36)"),
    @("4,961 (5.5)", "4,959 (5.5)"),
    @("18,247 (20)", "18,243 (20)"),
    @("22,703 (25)", "22,696 (25)"),
    @("23,634 (26)", "23,629 (26)"),
    @("20,713 (23)", "20,710 (23)"),
    @("64,281 (71)", "64,263 (71)"),
    @("22,781 (25)", "22,778 (25)"),
    @("16,297 (18)", "16,292 (18)"),
    @("30,604 (34)", "30,597 (34)"),
    @("43,357 (48)", "43,348 (48)"),
    @("74,546 (83)", "74,526 (83)"),
    @("14,991 (17)", "14,990 (17)"),
    @("76,549 (85)", "76,532 (85)"),
    @("13,285 (15)", "13,281 (15)"),
    @("29,461 (33)", "29,456 (33)"),
    @("54,052 (60)", "54,040 (60)"),
    @("6,745 (7.5)", "6,741 (7.5)"),
    @("15,370 (17)", "15,369 (17)"),
    @("66,477 (74)", "66,462 (74)"),
    @("8,411 (9.3)", "8,406 (9.3)"),
)

$successCount = 0
$failed = @()

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $result = $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)

    if ($result) {
        $successCount = $successCount + 1
    } else {
        $failed += $old
    }
}

Write-Output "Replacements applied: $successCount / $($replacements.Count)"
if ($failed.Count -gt 0) {
    Write-Output "Failed to find: $($failed -join '; ')"
}
